$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, shifting old rows 4-5 down to 5-6
$ws.Rows.Item(4).Insert()

# Update H3 text value
$ws.Range("H3").Value = "ORO BLANCO"

# Populate the newly inserted row 4 with data
$ws.Range("A4").Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Range("B4").Value = 45538
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("C4").Value = 45538
$ws.Range("C4").NumberFormat = "YYYY-MM-DD"
$ws.Range("D4").Value = 0.461
$ws.Range("E4").Value = 750000000
$ws.Range("F4").Value = 749769600
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "BNPDBC050924"
$ws.Range("I4").Value = "COMPRA"
$ws.Range("J4").Value = "RENTA VARIABLE"
